$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "first"
$ws.Range("C1").Value = "last"
$ws.Range("B1").Value = "name"

# Data rows
$ws.Range("A2").Value = "two"
$ws.Range("B2").Value = "pappu"
$ws.Range("C2").Value = "N"

$ws.Range("A3").Value = "two"
$ws.Range("B3").Value = "dabbu"
$ws.Range("C3").Value = "Y"

$ws.Range("A4").Value = "three"
$ws.Range("B4").Value = "dahakan"
$ws.Range("C4").Value = "NY"

$ws.Range("C3").Select()
